$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "French"
$ws.Range("D1").Value = "French"

$ws.Range("B2").Value = "comment"
$ws.Range("D2").Value = "quand"

$ws.Range("B3").Value = "quel/quelle"
$ws.Range("D3").Value = "qui"

$ws.Range("B4").Value = "où"
$ws.Range("D4").Value = "pourquoi"

$ws.Range("D1:D4").Select()
